# Applies:
#  1. Refresh the cached "datetimeFigureOut" date field text (4/2/2019 -> 4/9/2019)
#     on the slide master, every slide layout, and the notes master.
#  2. Reshape two "sequence lifeline" connectors on slide 1 (Straight Connector 10
#     and Straight Connector 89) to their new geometry / flip state.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: find the "Date" placeholder shape (ppPlaceholderDate = 16) inside a
# shape collection (slide master / custom layout / notes master all expose
# .Shapes the same way).
# ---------------------------------------------------------------------------
function Get-DatePlaceholderShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat -ne $null -and $shp.PlaceholderFormat.Type -eq 16) {
            return $shp
        }
    }
    return $null
}

function Set-DatePlaceholderText($shapes, [string]$newText) {
    $dateShape = Get-DatePlaceholderShape $shapes
    if ($dateShape -ne $null) {
        $dateShape.TextFrame.TextRange.Text = $newText
    }
}

$newDate = "4/9/2019"

# Slide master.
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every custom (slide) layout off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes $newDate
}

# ---------------------------------------------------------------------------
# Slide 1 connector geometry changes.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# "Straight Connector 10" -> flips horizontally and stretches to the new
# lifeline length.
$cxn10 = $slide.Shapes.Item("Straight Connector 10")
$cxn10.HorizontalFlip = -1
$cxn10.Left = 107.64142
$cxn10.Top = 117.304
$cxn10.Width = 0.5431
$cxn10.Height = 238.8913

# "Straight Connector 89" -> becomes perfectly vertical (0 width) and
# stretches to the same new lifeline length.
$cxn89 = $slide.Shapes.Item("Straight Connector 89")
$cxn89.Left = 211.4255
$cxn89.Top = 117.304
$cxn89.Width = 0.0
$cxn89.Height = 238.8913

Write-Output "Done."
